## Vscs sheet: insert a new row for "LOG prefix Name" before the existing
## "VSD FQDN" row (i.e. before current row 29), pushing all subsequent rows
## (and their comments, merged cells, data validations, etc.) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

## 1. Insert a blank row at position 29. Excel automatically shifts the
##    sheetData, mergeCells and dataValidations sqref ranges down by one,
##    but (in this engine) it does NOT relocate cell comments - those need
##    to be moved manually below.
$ws.Rows.Item(29).Insert()

## 2. Fill in the new row's label cell, matching the style already used by
##    its neighbouring rows (s="8" for col A, s="9" for cols B/C). Insert()
##    already copies the column-A style down correctly; B/C need a format
##    touch-up copied from the row below.
$ws.Range("A29").Value = "LOG prefix Name"
$ws.Range("B30:C30").Copy()
$ws.Range("B29:C29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

## 3. Relocate the comments that lived on rows 29-99 down to rows 30-100,
##    processing from the bottom up so that no comment is overwritten
##    before it has been read.
for ($r = 99; $r -ge 29; $r--) {
    $srcCell = $ws.Cells.Item($r, 1)
    $cmt = $srcCell.Comment
    if ($cmt -ne $null) {
        $txt = $cmt.Text()
        $cmt.Delete()
        $dstCell = $ws.Cells.Item($r + 1, 1)
        $dstCell.AddComment($txt)
    }
}

## 4. Add the new comment describing the new "LOG prefix Name" field.
$ws.Range("A29").AddComment("Log prefix name to pass in the vsc config")
